$d = $word.ActiveDocument

$replacements = @(
    @{old = "741×3=2223"; new = "576×2=1152"},
    @{old = "299×7=2093"; new = "175×6=1050"},
    @{old = "556×4=2224"; new = "842×7=5894"},
    @{old = "758×4=3032"; new = "514×4=2056"},
    @{old = "519×8=4152"; new = "311×2=622"},
    @{old = "633×4=2532"; new = "633×9=5697"},
    @{old = "941×6=5646"; new = "641×4=2564"},
    @{old = "541×3=1623"; new = "152×8=1216"},
    @{old = "716×2=1432"; new = "101×4=404"},
    @{old = "834×5=4170"; new = "984×7=6888"},
    @{old = "992×5=4960"; new = "441×6=2646"},
    @{old = "250×9=2250"; new = "897×4=3588"},
    @{old = "539×2=1078"; new = "849×6=5094"},
    @{old = "513×8=4104"; new = "676×8=5408"},
    @{old = "568×2=1136"; new = "380×6=2280"},
    @{old = "575×8=4600"; new = "928×9=8352"},
    @{old = "622×9=5598"; new = "442×6=2652"},
    @{old = "500×4=2000"; new = "739×6=4434"},
    @{old = "594×6=3564"; new = "170×8=1360"},
    @{old = "949×8=7592"; new = "102×9=918"},
    @{old = "350×5=1750"; new = "238×5=1190"},
    @{old = "937×7=6559"; new = "840×4=3360"},
    @{old = "966×7=6762"; new = "900×4=3600"},
    @{old = "425×5=2125"; new = "567×4=2268"},
    @{old = "857×2=1714"; new = "563×7=3941"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
